$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.091.06"
$ws.Range("E2").Value = "  +1.22%  "

# Row 3
$ws.Range("D3").Value = "'1.790.23"
$ws.Range("E3").Value = "  +1.56%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").Value = "'323.70"
$ws.Range("E5").Value = "  -1.09%  "

# Row 6
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.16%  "

# Row 7
$ws.Range("D7").Value = "'0.4282"
$ws.Range("E7").Value = "  -3.59%  "

# Row 8
$ws.Range("D8").Value = "'0.3617"
$ws.Range("E8").Value = "  -3.37%  "

# Row 9
$ws.Range("D9").Value = "'44.72"
$ws.Range("E9").Value = "  -1.48%  "

# Row 10
$ws.Range("D10").Value = "'0.07504"
$ws.Range("E10").Value = "  -3.64%  "

# Row 11
$ws.Range("D11").Value = "'1.111"
$ws.Range("E11").Value = "  -1.80%  "

# Row 12
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.03%  "

# Row 13
$ws.Range("D13").Value = "'21.59"
$ws.Range("E13").Value = "  -1.11%  "

# Row 14
$ws.Range("D14").Value = "'6.141"
$ws.Range("E14").Value = "  -1.19%  "

# Row 15
$ws.Range("D15").Value = "'7.314"
$ws.Range("E15").Value = "  -0.94%  "

# Row 16
$ws.Range("D16").Value = "'1.806.39"
$ws.Range("E16").Value = "  +2.68%  "

# Row 17
$ws.Range("D17").Value = "'91.98"
$ws.Range("E17").Value = "  +0.42%  "

# Row 18
$ws.Range("D18").Value = "'0.00001065"
$ws.Range("E18").Value = "  -1.88%  "

# Row 19
$ws.Range("D19").Value = "'0.06352"
$ws.Range("E19").Value = "  +1.57%  "

# Row 20
$ws.Range("E20").Value = "  -0.02%  "

# Row 21
$ws.Range("D21").Value = "'17.13"
$ws.Range("E21").Value = "  -1.58%  "

# Row 22
$ws.Range("D22").Value = "'5.969"
$ws.Range("E22").Value = "  -3.73%  "

# Row 23
$ws.Range("D23").Value = "'28.087.72"
$ws.Range("E23").Value = "  +1.11%  "

# Row 24
$ws.Range("D24").Value = "'11.35"
$ws.Range("E24").Value = "  -2.75%  "

# Row 25
$ws.Range("D25").Value = "'2.142"
$ws.Range("E25").Value = "  -8.22%  "

# Row 26
$ws.Range("D26").Value = "'158.91"
$ws.Range("E26").Value = "  +3.16%  "

# Row 27
$ws.Range("D27").Value = "'20.29"
$ws.Range("E27").Value = "  -2.90%  "

# Row 28
$ws.Range("D28").Value = "'2.009.35"
$ws.Range("E28").Value = "  +2.61%  "

# Row 29
$ws.Range("D29").Value = "'2.193"
$ws.Range("E29").Value = "  -7.79%  "

# Row 30
$ws.Range("D30").Value = "'127.20"
$ws.Range("E30").Value = "  -1.81%  "

# Row 31
$ws.Range("E31").Value = "  -4.05%  "

# Row 32
$ws.Range("D32").Value = "'5.792"
$ws.Range("E32").Value = "  -0.34%  "

# Row 33
$ws.Range("D33").Value = "'0.08980"
$ws.Range("E33").Value = "  -3.44%  "

# Row 34
$ws.Range("D34").Value = "'3.515"
$ws.Range("E34").Value = "  -4.20%  "

# Row 35
$ws.Range("D35").Value = "'12.63"
$ws.Range("E35").Value = "  -1.06%  "

# Row 36
$ws.Range("D36").Value = "'0.02342"
$ws.Range("E36").Value = "  +0.26%  "

# Row 37
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").Value = "'0.6473"
$ws.Range("E37").Value = "  -0.79%  "

# Row 38
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.071"
$ws.Range("E38").Value = "  -0.78%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2113"
$ws.Range("E39").Value = "  -3.73%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.06068"
$ws.Range("E40").Value = "  -1.67%  "

# Row 41
$ws.Range("D41").Value = "'1.183"
$ws.Range("E41").Value = "  -1.40%  "

# Row 42
$ws.Range("D42").Value = "'1.423"
$ws.Range("E42").Value = "  +0.43%  "

# Row 43
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = "  +0.15%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'7.908"
$ws.Range("E44").Value = "  -1.74%  "

# Row 45
$ws.Range("D45").Value = "'13.56"
$ws.Range("E45").Value = "  -2.24%  "

# Row 46
$ws.Range("D46").Value = "'0.5994"
$ws.Range("E46").Value = "  -0.80%  "

# Row 47
$ws.Range("D47").Value = "'3.711"
$ws.Range("E47").Value = "  -1.63%  "

# Row 48
$ws.Range("D48").Value = "'124.77"
$ws.Range("E48").Value = "  -1.37%  "

# Row 49
$ws.Range("D49").Value = "'1.983"
$ws.Range("E49").Value = "  -0.96%  "

# Row 50
$ws.Range("D50").Value = "'1.150"
$ws.Range("E50").Value = "  +0.08%  "

# Row 51
$ws.Range("D51").Value = "'0.06960"
$ws.Range("E51").Value = "  +0.60%  "
